$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 494, pushing existing rows 494:583 down to 495:584
$ws.Rows.Item(494).Insert()

# Populate the newly inserted row 494 with the new record
$ws.Range("A494").Value = 9
$ws.Range("B494").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C494").Value = "Metropolitana"
$ws.Range("D494").Value = 45015
$ws.Range("E494").Value = 13
$ws.Range("F494").Value = 100112012
$ws.Range("G494").Value = "Espinaca"
$ws.Range("H494").Value = "Sin especificar"
$ws.Range("I494").Value = "Primera"
$ws.Range("J494").Value = 160
$ws.Range("K494").Value = 7000
$ws.Range("L494").Value = 8000
$ws.Range("M494").Value = 7500
$ws.Range("N494").Value = "$/cuna 10 kilos"
$ws.Range("O494").Value = "Provincia de Chacabuco"
$ws.Range("P494").Value = 750
$ws.Range("Q494").Value = 10
$ws.Range("R494").Value = "Hortaliza"

# Match the date number format used by the rest of column D
$ws.Range("D494").NumberFormat = $ws.Range("D495").NumberFormat
